$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source added one new daily-price record for this market/variety.
# In the canonical sheet it lands at row 57 (sorted position), pushing the
# existing rows 57:140 down to 58:141 (dimension grows from T140 to T141).
$ws.Rows("57").Insert()

# Populate the newly inserted row 57 with the new record.
$ws.Cells.Item(57, 1).Value = 11
$ws.Cells.Item(57, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(57, 3).Value = 'Bíobío'
$ws.Cells.Item(57, 4).Value = 44705
$ws.Cells.Item(57, 5).Value = 8
$ws.Cells.Item(57, 6).Value = 'Fruta'
$ws.Cells.Item(57, 7).Value = 100101
$ws.Cells.Item(57, 8).Value = 'Berries'
$ws.Cells.Item(57, 9).Value = 100101007
$ws.Cells.Item(57, 10).Value = 'Kiwi'
$ws.Cells.Item(57, 11).Value = 'Hayward'
$ws.Cells.Item(57, 12).Value = 'Primera'
$ws.Cells.Item(57, 13).Value = 220
$ws.Cells.Item(57, 14).Value = 10000
$ws.Cells.Item(57, 15).Value = 11000
$ws.Cells.Item(57, 16).Value = 10455
$ws.Cells.Item(57, 17).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(57, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(57, 19).Value = 581
$ws.Cells.Item(57, 20).Value = 18
